$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh (ref -> new text value).
# Kept as literal strings: most "prices" (e.g. "66.308.14") are not valid
# Excel numbers (thousands separated by ".") and must round-trip byte-for-byte,
# so every cell is forced to Text format before the write and reset to the
# Normal style afterwards (matches the unstyled inlineStr cells in the source).
$updates = @(
    @{ Ref = "D2"; Value = '66.308.14' },
    @{ Ref = "E2"; Value = '  -5.94%  ' },
    @{ Ref = "D3"; Value = '3.179.83' },
    @{ Ref = "E3"; Value = '  -9.63%  ' },
    @{ Ref = "D4"; Value = '0.999' },
    @{ Ref = "E4"; Value = '  -0.04%  ' },
    @{ Ref = "D5"; Value = '571.42' },
    @{ Ref = "E5"; Value = '  -6.48%  ' },
    @{ Ref = "D6"; Value = '147.97' },
    @{ Ref = "E6"; Value = '  -14.61%  ' },
    @{ Ref = "D7"; Value = '0.999' },
    @{ Ref = "E7"; Value = '  -0.05%  ' },
    @{ Ref = "D8"; Value = '3.171.72' },
    @{ Ref = "E8"; Value = '  -9.67%  ' },
    @{ Ref = "D9"; Value = '0.537' },
    @{ Ref = "E9"; Value = '  -11.88%  ' },
    @{ Ref = "D10"; Value = '0.167' },
    @{ Ref = "E10"; Value = '  -14.98%  ' },
    @{ Ref = "D11"; Value = '6.27' },
    @{ Ref = "E11"; Value = '  -12.29%  ' },
    @{ Ref = "D12"; Value = '0.490' },
    @{ Ref = "E12"; Value = '  -16.94%  ' },
    @{ Ref = "D13"; Value = '37.94' },
    @{ Ref = "E13"; Value = '  -18.37%  ' },
    @{ Ref = "D14"; Value = '0.0000239' },
    @{ Ref = "E14"; Value = '  -13.59%  ' },
    @{ Ref = "D15"; Value = '3.687.18' },
    @{ Ref = "E15"; Value = '  -9.76%  ' },
    @{ Ref = "D16"; Value = '66.273.36' },
    @{ Ref = "E16"; Value = '  -5.94%  ' },
    @{ Ref = "D17"; Value = '3.174.14' },
    @{ Ref = "E17"; Value = '  -9.65%  ' },
    @{ Ref = "E18"; Value = '  -6.89%  ' },
    @{ Ref = "D19"; Value = '528.54' },
    @{ Ref = "E19"; Value = '  -13.92%  ' },
    @{ Ref = "E20"; Value = '  -16.90%  ' },
    @{ Ref = "D21"; Value = '14.84' },
    @{ Ref = "E21"; Value = '  -16.63%  ' },
    @{ Ref = "E22"; Value = '  -15.34%  ' },
    @{ Ref = "D23"; Value = '7.68' },
    @{ Ref = "E23"; Value = '  -14.87%  ' },
    @{ Ref = "D24"; Value = '84.25' },
    @{ Ref = "E24"; Value = '  -14.84%  ' },
    @{ Ref = "D25"; Value = '13.25' },
    @{ Ref = "E25"; Value = '  -15.62%  ' },
    @{ Ref = "D26"; Value = '0.999' },
    @{ Ref = "E26"; Value = '  -0.09%  ' },
    @{ Ref = "D27"; Value = '3.09' },
    @{ Ref = "E27"; Value = '  -18.05%  ' },
    @{ Ref = "D28"; Value = '2.14' },
    @{ Ref = "E28"; Value = '  -17.32%  ' },
    @{ Ref = "D29"; Value = '7.94' },
    @{ Ref = "E29"; Value = '  -13.79%  ' },
    @{ Ref = "E30"; Value = '  -15.11%  ' },
    @{ Ref = "D31"; Value = '2.55' },
    @{ Ref = "E31"; Value = '  -15.29%  ' },
    @{ Ref = "E32"; Value = '  -14.81%  ' },
    @{ Ref = "E33"; Value = '  -21.23%  ' },
    @{ Ref = "B34"; Value = 'Bittensor' },
    @{ Ref = "C34"; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' },
    @{ Ref = "D34"; Value = '519.35' },
    @{ Ref = "E34"; Value = '  -15.66%  ' },
    @{ Ref = "B35"; Value = 'NEARProtocol' },
    @{ Ref = "C35"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Ref = "D35"; Value = '5.58' },
    @{ Ref = "E35"; Value = '  -18.67%  ' },
    @{ Ref = "B36"; Value = 'FirstDigitalUSD' },
    @{ Ref = "C36"; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' },
    @{ Ref = "D36"; Value = '1.00' },
    @{ Ref = "E36"; Value = '  +0.19%  ' },
    @{ Ref = "D37"; Value = '52.99' },
    @{ Ref = "E37"; Value = '  -6.98%  ' },
    @{ Ref = "D38"; Value = '0.0847' },
    @{ Ref = "E38"; Value = '  -15.81%  ' },
    @{ Ref = "D39"; Value = '9.01' },
    @{ Ref = "E39"; Value = '  -16.80%  ' },
    @{ Ref = "D40"; Value = '0.0407' },
    @{ Ref = "E40"; Value = '  -18.15%  ' },
    @{ Ref = "D41"; Value = '0.123' },
    @{ Ref = "E41"; Value = '  -14.79%  ' },
    @{ Ref = "D42"; Value = '2.868.02' },
    @{ Ref = "E42"; Value = '  -14.99%  ' },
    @{ Ref = "D43"; Value = '2.57' },
    @{ Ref = "E43"; Value = '  -26.41%  ' },
    @{ Ref = "D44"; Value = '0.257' },
    @{ Ref = "E44"; Value = '  -17.43%  ' },
    @{ Ref = "D45"; Value = '0.0₃0574' },
    @{ Ref = "E45"; Value = '  -22.87%  ' },
    @{ Ref = "D47"; Value = '25.75' },
    @{ Ref = "E47"; Value = '  -20.03%  ' },
    @{ Ref = "E48"; Value = '  -20.85%  ' },
    @{ Ref = "E49"; Value = '  -19.66%  ' },
    @{ Ref = "E50"; Value = '  -14.09%  ' },
    @{ Ref = "D51"; Value = '122.38' },
    @{ Ref = "E51"; Value = '  -8.23%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
